{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph,\n// the copyright/footer paragraph that follows it, and the blank paragraph that\n// immediately precedes them (right after the \"LOB1019...\" requirement line).\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two text paragraphs that must be deleted outright.\nconst viewIdx = items.findIndex(p => p.text.trim() === \"Ver no Jupiter Salvar em pdf Salvar em docx\");\nconst copyrightIdx = items.findIndex(p => p.text.trim().startsWith(\"\u00a9 2020\"));\n\nconst toDelete = [];\nif (viewIdx !== -1) toDelete.push(items[viewIdx]);\nif (copyrightIdx !== -1) toDelete.push(items[copyrightIdx]);\n\n// The empty paragraph right before the \"Ver no Jupiter...\" paragraph (if blank)\n// also needs to go, so the requirement line is followed directly by the\n// pre-existing trailing blank paragraph / page break, matching the diff.\nif (viewIdx > 0) {\n    const prev = items[viewIdx - 1];\n    if (prev.text.trim() === \"\") {\n        toDelete.push(prev);\n    }\n}\n\nfor (const p of toDelete) {\n    p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph,\n# the copyright/footer paragraph that follows it, and the blank paragraph that\n# immediately precedes them (right after the \"LOB1019...\" requirement line).\n\n$d = $word.ActiveDocument\n\n$viewIndex = -1\n$copyrightIndex = -1\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t.Contains(\"Ver no Jupiter Salvar em pdf Salvar em docx\")) {\n        $viewIndex = $i\n    }\n    if ($t.Contains(\"Contact: luizeleno@usp.br\")) {\n        $copyrightIndex = $i\n    }\n}\n\n# Delete from the highest index down so earlier indices stay valid.\nif ($copyrightIndex -ge 1) {\n    $d.Paragraphs.Item($copyrightIndex).Range.Delete()\n}\nif ($viewIndex -ge 1) {\n    $d.Paragraphs.Item($viewIndex).Range.Delete()\n}\n\n# The blank paragraph that used to sit right before the \"Ver no Jupiter...\"\n# paragraph is also removed by the edit.\nif ($viewIndex -ge 2) {\n    $blankIndex = $viewIndex - 1\n    $blankText = $d.Paragraphs.Item($blankIndex).Range.Text.Trim()\n    if ($blankText -eq \"\") {\n        $d.Paragraphs.Item($blankIndex).Range.Delete()\n    }\n}\n"}
